# Auto-generated edit script: update cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.591.32"
$ws.Range("E2").Value = "  -4.36%  "

$ws.Range("D3").Value = "'3.259.35"
$ws.Range("E3").Value = "  -5.33%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'556.77"
$ws.Range("E5").Value = "  -2.76%  "

$ws.Range("D6").Value = "'181.17"
$ws.Range("E6").Value = "  -4.28%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").Value = "'3.255.04"
$ws.Range("E9").Value = "  -5.13%  "

$ws.Range("E10").Value = "  -8.14%  "

$ws.Range("E11").Value = "  -4.48%  "

$ws.Range("D12").Value = "'47.30"
$ws.Range("E12").Value = "  -6.85%  "

$ws.Range("E13").Value = "  -6.24%  "

$ws.Range("D14").Value = "'637.46"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("E15").Value = "  -5.34%  "

$ws.Range("D16").Value = "'3.788.22"
$ws.Range("E16").Value = "  -4.92%  "

$ws.Range("D17").Value = "'65.533.82"
$ws.Range("E17").Value = "  -4.14%  "

$ws.Range("E18").Value = "  -3.18%  "

$ws.Range("D19").Value = "'17.69"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'3.264.64"
$ws.Range("E20").Value = "  -5.28%  "

$ws.Range("D21").Value = "'11.35"
$ws.Range("E21").Value = "  -7.04%  "

$ws.Range("E22").Value = "  -3.53%  "

$ws.Range("D23").Value = "'17.83"
$ws.Range("E23").Value = "  +0.56%  "

$ws.Range("D24").Value = "'105.61"
$ws.Range("E24").Value = "  +6.90%  "

$ws.Range("D25").Value = "'4.96"
$ws.Range("E25").Value = "  -6.78%  "

$ws.Range("D26").Value = "'3.98"
$ws.Range("E26").Value = "  -6.26%  "

$ws.Range("E27").Value = "  -5.52%  "

$ws.Range("D28").Value = "'9.51"
$ws.Range("E28").Value = "  -2.78%  "

$ws.Range("E29").Value = "  -4.71%  "

$ws.Range("D30").Value = "'30.29"
$ws.Range("E30").Value = "  -5.77%  "

$ws.Range("D31").Value = "'4.05"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("D32").Value = "'6.31"
$ws.Range("E32").Value = "  -5.25%  "

$ws.Range("D33").Value = "'11.01"
$ws.Range("E33").Value = "  -4.21%  "

$ws.Range("D34").Value = "'552.69"
$ws.Range("E34").Value = "  +10.71%  "

$ws.Range("E35").Value = "  -2.57%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'56.96"
$ws.Range("E37").Value = "  -6.13%  "

$ws.Range("D38").Value = "'3.599.88"
$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = "  +6.01%  "

$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  -1.73%  "

$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "'0.0₃0720"
$ws.Range("E41").Value = "  -7.56%  "

$ws.Range("D42").Value = "'2.73"
$ws.Range("E42").Value = "  -5.14%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.129"
$ws.Range("E43").Value = "  -2.01%  "

$ws.Range("D44").Value = "'31.98"
$ws.Range("E44").Value = "  -6.19%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.32"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.335"
$ws.Range("E46").Value = "  -7.94%  "

$ws.Range("D47").Value = "'0.0415"
$ws.Range("E47").Value = "  -4.39%  "

$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  -6.33%  "

$ws.Range("E49").Value = "  -3.09%  "

$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("D51").Value = "'1.24"
$ws.Range("E51").Value = "  +2.16%  "
